$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", matching the bold/bordered header style used by B1:E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Fill F2:F144 with the recorded time_taken timestamps
$arr = New-Object 'object[,]' 143,1
$arr[0,0] = "2021-10-05 13:41:43.864191"
$arr[1,0] = "2021-10-05 13:41:43.864204"
$arr[2,0] = "2021-10-05 13:41:43.864207"
$arr[3,0] = "2021-10-05 13:41:43.864210"
$arr[4,0] = "2021-10-05 13:41:43.864212"
$arr[5,0] = "2021-10-05 13:41:43.864215"
$arr[6,0] = "2021-10-05 13:41:43.864218"
$arr[7,0] = "2021-10-05 13:41:43.864224"
$arr[8,0] = "2021-10-05 13:41:43.864227"
$arr[9,0] = "2021-10-05 13:41:43.864229"
$arr[10,0] = "2021-10-05 13:41:43.864232"
$arr[11,0] = "2021-10-05 13:41:43.864234"
$arr[12,0] = "2021-10-05 13:41:43.864237"
$arr[13,0] = "2021-10-05 13:41:43.864239"
$arr[14,0] = "2021-10-05 13:41:43.864242"
$arr[15,0] = "2021-10-05 13:41:43.864245"
$arr[16,0] = "2021-10-05 13:41:43.864247"
$arr[17,0] = "2021-10-05 13:41:43.864250"
$arr[18,0] = "2021-10-05 13:41:43.864252"
$arr[19,0] = "2021-10-05 13:41:43.864255"
$arr[20,0] = "2021-10-05 13:41:43.864257"
$arr[21,0] = "2021-10-05 13:41:43.864260"
$arr[22,0] = "2021-10-05 13:41:43.864264"
$arr[23,0] = "2021-10-05 13:41:43.864266"
$arr[24,0] = "2021-10-05 13:41:43.864269"
$arr[25,0] = "2021-10-05 13:41:43.864272"
$arr[26,0] = "2021-10-05 13:41:43.864275"
$arr[27,0] = "2021-10-05 13:41:43.864277"
$arr[28,0] = "2021-10-05 13:41:43.864280"
$arr[29,0] = "2021-10-05 13:41:43.864283"
$arr[30,0] = "2021-10-05 13:41:43.864285"
$arr[31,0] = "2021-10-05 13:41:43.864288"
$arr[32,0] = "2021-10-05 13:41:43.864291"
$arr[33,0] = "2021-10-05 13:41:43.864293"
$arr[34,0] = "2021-10-05 13:41:43.864296"
$arr[35,0] = "2021-10-05 13:41:43.864298"
$arr[36,0] = "2021-10-05 13:41:43.864301"
$arr[37,0] = "2021-10-05 13:41:43.864303"
$arr[38,0] = "2021-10-05 13:41:43.864306"
$arr[39,0] = "2021-10-05 13:41:43.864308"
$arr[40,0] = "2021-10-05 13:41:43.864311"
$arr[41,0] = "2021-10-05 13:41:43.864314"
$arr[42,0] = "2021-10-05 13:41:43.864317"
$arr[43,0] = "2021-10-05 13:41:43.864319"
$arr[44,0] = "2021-10-05 13:41:43.864321"
$arr[45,0] = "2021-10-05 13:41:43.864324"
$arr[46,0] = "2021-10-05 13:41:43.864326"
$arr[47,0] = "2021-10-05 13:41:43.864329"
$arr[48,0] = "2021-10-05 13:41:43.864331"
$arr[49,0] = "2021-10-05 13:41:43.864334"
$arr[50,0] = "2021-10-05 13:41:43.864336"
$arr[51,0] = "2021-10-05 13:41:43.864339"
$arr[52,0] = "2021-10-05 13:41:43.864342"
$arr[53,0] = "2021-10-05 13:41:43.864345"
$arr[54,0] = "2021-10-05 13:41:43.864347"
$arr[55,0] = "2021-10-05 13:41:43.864350"
$arr[56,0] = "2021-10-05 13:41:43.864352"
$arr[57,0] = "2021-10-05 13:41:43.864355"
$arr[58,0] = "2021-10-05 13:41:43.864357"
$arr[59,0] = "2021-10-05 13:41:43.864360"
$arr[60,0] = "2021-10-05 13:41:43.864362"
$arr[61,0] = "2021-10-05 13:41:43.864365"
$arr[62,0] = "2021-10-05 13:41:43.864368"
$arr[63,0] = "2021-10-05 13:41:43.864370"
$arr[64,0] = "2021-10-05 13:41:43.864374"
$arr[65,0] = "2021-10-05 13:41:43.864377"
$arr[66,0] = "2021-10-05 13:41:43.864379"
$arr[67,0] = "2021-10-05 13:41:43.864382"
$arr[68,0] = "2021-10-05 13:41:43.864384"
$arr[69,0] = "2021-10-05 13:41:43.864387"
$arr[70,0] = "2021-10-05 13:41:43.864389"
$arr[71,0] = "2021-10-05 13:41:43.864392"
$arr[72,0] = "2021-10-05 13:41:43.864394"
$arr[73,0] = "2021-10-05 13:41:43.864397"
$arr[74,0] = "2021-10-05 13:41:43.864399"
$arr[75,0] = "2021-10-05 13:41:43.864402"
$arr[76,0] = "2021-10-05 13:41:43.864407"
$arr[77,0] = "2021-10-05 13:41:43.864410"
$arr[78,0] = "2021-10-05 13:41:43.864413"
$arr[79,0] = "2021-10-05 13:41:43.864415"
$arr[80,0] = "2021-10-05 13:41:43.864418"
$arr[81,0] = "2021-10-05 13:41:43.864420"
$arr[82,0] = "2021-10-05 13:41:43.864423"
$arr[83,0] = "2021-10-05 13:41:43.864426"
$arr[84,0] = "2021-10-05 13:41:43.864428"
$arr[85,0] = "2021-10-05 13:41:43.864431"
$arr[86,0] = "2021-10-05 13:41:43.864434"
$arr[87,0] = "2021-10-05 13:41:43.864436"
$arr[88,0] = "2021-10-05 13:41:43.864439"
$arr[89,0] = "2021-10-05 13:41:43.864441"
$arr[90,0] = "2021-10-05 13:41:43.864444"
$arr[91,0] = "2021-10-05 13:41:43.864446"
$arr[92,0] = "2021-10-05 13:41:43.864450"
$arr[93,0] = "2021-10-05 13:41:43.864453"
$arr[94,0] = "2021-10-05 13:41:43.864456"
$arr[95,0] = "2021-10-05 13:41:43.864458"
$arr[96,0] = "2021-10-05 13:41:43.864461"
$arr[97,0] = "2021-10-05 13:41:43.864463"
$arr[98,0] = "2021-10-05 13:41:43.864466"
$arr[99,0] = "2021-10-05 13:41:43.864468"
$arr[100,0] = "2021-10-05 13:41:43.864471"
$arr[101,0] = "2021-10-05 13:41:43.864473"
$arr[102,0] = "2021-10-05 13:41:43.864476"
$arr[103,0] = "2021-10-05 13:41:43.864479"
$arr[104,0] = "2021-10-05 13:41:43.864481"
$arr[105,0] = "2021-10-05 13:41:43.864484"
$arr[106,0] = "2021-10-05 13:41:43.864486"
$arr[107,0] = "2021-10-05 13:41:43.864489"
$arr[108,0] = "2021-10-05 13:41:43.864493"
$arr[109,0] = "2021-10-05 13:41:43.864496"
$arr[110,0] = "2021-10-05 13:41:43.864499"
$arr[111,0] = "2021-10-05 13:41:43.864502"
$arr[112,0] = "2021-10-05 13:41:43.864504"
$arr[113,0] = "2021-10-05 13:41:43.864507"
$arr[114,0] = "2021-10-05 13:41:43.864509"
$arr[115,0] = "2021-10-05 13:41:43.864512"
$arr[116,0] = "2021-10-05 13:41:43.864514"
$arr[117,0] = "2021-10-05 13:41:43.864517"
$arr[118,0] = "2021-10-05 13:41:43.864519"
$arr[119,0] = "2021-10-05 13:41:43.864522"
$arr[120,0] = "2021-10-05 13:41:43.864524"
$arr[121,0] = "2021-10-05 13:41:43.864528"
$arr[122,0] = "2021-10-05 13:41:43.864530"
$arr[123,0] = "2021-10-05 13:41:43.864533"
$arr[124,0] = "2021-10-05 13:41:43.864535"
$arr[125,0] = "2021-10-05 13:41:43.864538"
$arr[126,0] = "2021-10-05 13:41:43.864540"
$arr[127,0] = "2021-10-05 13:41:43.864543"
$arr[128,0] = "2021-10-05 13:41:43.864547"
$arr[129,0] = "2021-10-05 13:41:43.864550"
$arr[130,0] = "2021-10-05 13:41:43.864553"
$arr[131,0] = "2021-10-05 13:41:43.864555"
$arr[132,0] = "2021-10-05 13:41:43.864558"
$arr[133,0] = "2021-10-05 13:41:43.864561"
$arr[134,0] = "2021-10-05 13:41:43.864563"
$arr[135,0] = "2021-10-05 13:41:43.864566"
$arr[136,0] = "2021-10-05 13:41:43.864568"
$arr[137,0] = "2021-10-05 13:41:43.864571"
$arr[138,0] = "2021-10-05 13:41:43.864573"
$arr[139,0] = "2021-10-05 13:41:43.864576"
$arr[140,0] = "2021-10-05 13:41:43.864579"
$arr[141,0] = "2021-10-05 13:41:43.864581"
$arr[142,0] = "2021-10-05 13:41:43.864584"
$ws.Range("F2:F144").Value = $arr

"Done: F1=" + $ws.Range("F1").Value2 + " F144=" + $ws.Range("F144").Value2
